$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin rows: Price (D) and Volume 1h (E) values, with a couple of
# rows (31/32 and 38/39) reordered entirely (Coin name + Link + Price + Volume).
# Leading "'" forces text entry (these prices/percentages are text, not numbers);
# resetting the Style back to Normal afterwards drops the quote-prefix/text
# formatting that Excel applies automatically so the cell keeps its original style.

$ws.Range("D2").Value = "'29.414.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.07%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.842.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'239.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.31%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6259"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.52%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D8").Value = "'0.07396"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.89%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2892"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.52%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'24.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.54%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07717"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.28%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.835.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.36%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.967"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.65%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6716"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.00001035"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.12%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'81.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.38%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'6.273"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.28%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'29.380.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.10%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'234.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.58%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.31%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D22").Value = "'7.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.84%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.72%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'156.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.40%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'8.479"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.28%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.1346"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.44%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'17.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.08%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.07267"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +12.11%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.497"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.56%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.478"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.42%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'Filecoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'4.032"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.43%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'4.038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.46%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.163"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.89%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.816"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.05%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.7100"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.04%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'2.582"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.04%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.01838"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.48%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'MXToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'2.787"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.73%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'Maker"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.233.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.44%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'6.804"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.42%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.9550"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.41%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D43").Value = "'1.992.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.64%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'101.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.14%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'65.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.97%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000119"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.84%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.700"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.03%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'6.957"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.74%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'8.935"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.17%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.1133"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.16%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.3884"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.56%  "
$ws.Range("E51").Style = "Normal"
